$d = $word.ActiveDocument

$pairs = @(
    @("31×29=", "99×62="),
    @("35×68=", "75×97="),
    @("77×12=", "46×15="),
    @("77×33=", "36×47="),
    @("25×12=", "72×59="),
    @("96×52=", "30×72="),
    @("15×51=", "13×30="),
    @("85×53=", "73×53="),
    @("34×35=", "86×29="),
    @("74×29=", "92×72="),
    @("70×62=", "19×26="),
    @("24×80=", "29×26="),
    @("54×51=", "17×31="),
    @("43×85=", "55×34="),
    @("19×42=", "75×31="),
    @("65×89=", "68×18="),
    @("65×34=", "68×56="),
    @("87×49=", "52×72="),
    @("46×54=", "57×45="),
    @("66×82=", "75×16="),
    @("53×15=", "22×11="),
    @("43×15=", "59×63="),
    @("44×15=", "33×49="),
    @("20×66=", "71×45="),
    @("27×42=", "21×73=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
